$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.349.83"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "2.927.57"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.501"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").Value = "3.411.43"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "61.339.60"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "2.931.70"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "431.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -5.07%  "
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").Value = "0.0₃0878"
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.123"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.42%  "
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("D44").Value = "2.697.71"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "366.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.29%  "
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.125"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.87%  "
